$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.651.14"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.923.15"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.85"
$ws.Range("E5").Value = "  +4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4055"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08100"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.004"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.40"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.964.90"
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.002"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.201"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.14"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06841"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001029"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.56"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.014"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "29.676.87"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.576"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.79"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.161"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "2.193.56"
$ws.Range("E25").Value = "  +1.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.588"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.00"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.89"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.067"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.40"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.002"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09597"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.523"
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.397"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.541"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06545"
$ws.Range("E36").Value = "  +7.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02264"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.201"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5911"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.65"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.886"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1833"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.483"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.274"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.28"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07472"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5502"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.44"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.400"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.97"
$ws.Range("E51").Value = "  -0.48%  "

$ws.Range("D2:E51").ClearFormats()
